# Updated symbol list refresh: write new Price (D) and Volume(1h) (E) text values.
# Leading apostrophe forces literal text entry so formatted numeric strings
# (trailing zeros, percent signs, etc.) are preserved exactly, matching the
# scraped-text storage already used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'275.43"
$ws.Range("E2").Value = "'-1.66%"
# Row 3
$ws.Range("D3").Value = "'27.22"
$ws.Range("E3").Value = "'0.97%"
# Row 4
$ws.Range("D4").Value = "'4.762"
$ws.Range("E4").Value = "'-3.74%"
# Row 5
$ws.Range("D5").Value = "'0.06332"
$ws.Range("E5").Value = "'-1.28%"
# Row 6
$ws.Range("D6").Value = "'6.936"
$ws.Range("E6").Value = "'-0.65%"
# Row 7
$ws.Range("D7").Value = "'1.340"
$ws.Range("E7").Value = "'30.43%"
# Row 8
$ws.Range("D8").Value = "'0.8779"
$ws.Range("E8").Value = "'-1.27%"
# Row 9
$ws.Range("D9").Value = "'0.1508"
$ws.Range("E9").Value = "'1.02%"
# Row 10
$ws.Range("D10").Value = "'0.05043"
$ws.Range("E10").Value = "'-2.82%"
# Row 11
$ws.Range("D11").Value = "'0.07586"
$ws.Range("E11").Value = "'2.41%"
# Row 12
$ws.Range("D12").Value = "'0.02934"
$ws.Range("E12").Value = "'-6.81%"
# Row 13
$ws.Range("D13").Value = "'0.09004"
$ws.Range("E13").Value = "'-0.81%"
# Row 14
$ws.Range("D14").Value = "'0.001560"
$ws.Range("E14").Value = "'0.62%"
# Row 15
$ws.Range("D15").Value = "'0.0006392"
$ws.Range("E15").Value = "'1.25%"
# Row 16
$ws.Range("D16").Value = "'0.005853"
$ws.Range("E16").Value = "'-3.33%"
# Row 17
$ws.Range("D17").Value = "'3.448"
$ws.Range("E17").Value = "'-1.36%"
# Row 18
$ws.Range("D18").Value = "'3.293"
$ws.Range("E18").Value = "'-1.69%"
# Row 21
$ws.Range("E21").Value = "'0.83%"
# Row 22
$ws.Range("D22").Value = "'3.922"
$ws.Range("E22").Value = "'-0.40%"
# Row 23
$ws.Range("D23").Value = "'0.04414"
$ws.Range("E23").Value = "'1.01%"
# Row 24
$ws.Range("D24").Value = "'0.001171"
$ws.Range("E24").Value = "'-1.09%"
# Row 25
$ws.Range("D25").Value = "'0.003856"
$ws.Range("E25").Value = "'4.50%"
# Row 26
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("E26").Value = "'-0.47%"
# Row 27
$ws.Range("E27").Value = "'13.63%"
# Row 40
$ws.Range("D40").Value = "'0.04117"
$ws.Range("E40").Value = "'0.44%"
# Row 41
$ws.Range("D41").Value = "'0.006806"
$ws.Range("E41").Value = "'2.05%"
# Row 42
$ws.Range("D42").Value = "'0.1176"
$ws.Range("E42").Value = "'-0.41%"
# Row 43
$ws.Range("D43").Value = "'0.002087"
$ws.Range("E43").Value = "'-11.86%"
# Row 44
$ws.Range("D44").Value = "'0.01154"
$ws.Range("E44").Value = "'-8.18%"
# Row 45
$ws.Range("D45").Value = "'0.00005162"
$ws.Range("E45").Value = "'-2.16%"
# Row 46
$ws.Range("D46").Value = "'1.486"
$ws.Range("E46").Value = "'-36.89%"
# Row 47
$ws.Range("D47").Value = "'0.02294"
$ws.Range("E47").Value = "'2.08%"
